$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 2737.125
$ws.Range("I40").Value = 2047.2727
$ws.Range("K40").Value = 2047.2727
$ws.Range("M40").Value = -1872.2727
$ws.Range("H53").Value = 296.6875
$ws.Range("I53").Value = 156.9
$ws.Range("J53").Value = 529.6667
$ws.Range("K53").Value = 156.9
$ws.Range("L53").Value = 529.6667
$ws.Range("M53").Value = 480.1
$ws.Range("N53").Value = -1803.6667
$ws.Range("H62").Value = 6835.2856
$ws.Range("I62").Value = 6172
$ws.Range("J62").Value = 7498.5713
$ws.Range("K62").Value = 6172
$ws.Range("L62").Value = 7498.5713
$ws.Range("M62").Value = -5548
$ws.Range("N62").Value = -8746.5713
$ws.Range("H65").Value = 6835.2856
$ws.Range("I65").Value = 6172
$ws.Range("J65").Value = 7498.5713
$ws.Range("K65").Value = 30860
$ws.Range("L65").Value = 37492.85649999999
$ws.Range("M65").Value = -27740
$ws.Range("N65").Value = -43732.85649999999
$ws.Range("H74").Value = 4536.7915
$ws.Range("I74").Value = 3625.7273
$ws.Range("K74").Value = 3625.7273
$ws.Range("M74").Value = -2689.7273
$ws.Range("H76").Value = 11010.85
$ws.Range("I76").Value = 18846.857
$ws.Range("K76").Value = 18846.857
$ws.Range("M76").Value = -18531.857
$ws.Range("H77").Value = 4536.7915
$ws.Range("I77").Value = 3625.7273
$ws.Range("K77").Value = 18128.6365
$ws.Range("M77").Value = -13448.6365
$ws.Range("H79").Value = 11010.85
$ws.Range("I79").Value = 18846.857
$ws.Range("K79").Value = 18846.857
$ws.Range("M79").Value = -17754.857
$ws.Range("H113").Value = 8639.799999999999
$ws.Range("I113").Value = 6066.3335
$ws.Range("K113").Value = 6066.3335
$ws.Range("M113").Value = -2812.3335
$ws.Range("H137").Value = 3268.6667
$ws.Range("J137").Value = 3316.6667
$ws.Range("L137").Value = 9950.000100000001
$ws.Range("N137").Value = -15050.0001
$ws.Range("H138").Value = 5468852.5
$ws.Range("J138").Value = 7251380.5
$ws.Range("L138").Value = 21754141.5
$ws.Range("N138").Value = -21764421.5

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 40004284
$ws.Range("I61").Value = 66669576
$ws.Range("J61").Value = 6345
$ws.Range("K61").Value = 66669576
$ws.Range("L61").Value = 6345
$ws.Range("M61").Value = -66669364
$ws.Range("N61").Value = -6769
$ws.Range("H74").Value = 32261158
$ws.Range("I74").Value = 71430184
$ws.Range("J74").Value = 4312.353
$ws.Range("K74").Value = 71430184
$ws.Range("L74").Value = 4312.353
$ws.Range("M74").Value = -71429310
$ws.Range("N74").Value = -6060.353
$ws.Range("H77").Value = 32261158
$ws.Range("I77").Value = 71430184
$ws.Range("J77").Value = 4312.353
$ws.Range("K77").Value = 357150920
$ws.Range("L77").Value = 21561.765
$ws.Range("M77").Value = -357146552
$ws.Range("N77").Value = -30297.765
$ws.Range("H122").Value = 3485.4167
$ws.Range("I122").Value = 1788.1
$ws.Range("J122").Value = 4697.7856
$ws.Range("K122").Value = 5364.299999999999
$ws.Range("L122").Value = 14093.3568
$ws.Range("M122").Value = -2914.299999999999
$ws.Range("N122").Value = -18993.3568
$ws.Range("H136").Value = 40004284
$ws.Range("I136").Value = 66669576
$ws.Range("J136").Value = 6345
$ws.Range("K136").Value = 200008728
$ws.Range("L136").Value = 19035
$ws.Range("M136").Value = -200006178
$ws.Range("N136").Value = -24135

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 14091.36
$ws.Range("I86").Value = 7145.5264
$ws.Range("J86").Value = 36086.5
$ws.Range("K86").Value = 7145.5264
$ws.Range("L86").Value = 36086.5
$ws.Range("M86").Value = -6022.5264
$ws.Range("N86").Value = -38332.5
$ws.Range("H89").Value = 14091.36
$ws.Range("I89").Value = 7145.5264
$ws.Range("J89").Value = 36086.5
$ws.Range("K89").Value = 35727.632
$ws.Range("L89").Value = 180432.5
$ws.Range("M89").Value = -30111.632
$ws.Range("N89").Value = -191664.5
$ws.Range("H99").Value = 3137.963
$ws.Range("I99").Value = 2298.0833
$ws.Range("K99").Value = 2298.0833
$ws.Range("M99").Value = -800.0832999999998
$ws.Range("H134").Value = 5471.591
$ws.Range("I134").Value = 5241.7144
$ws.Range("J134").Value = 5873.875
$ws.Range("K134").Value = 15725.1432
$ws.Range("L134").Value = 17621.625
$ws.Range("M134").Value = -13190.1432
$ws.Range("N134").Value = -22691.625

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 7333.8
$ws.Range("I99").Value = 7955.4546
$ws.Range("K99").Value = 7955.4546
$ws.Range("M99").Value = -6457.4546
$ws.Range("H105").Value = 7641.1665
$ws.Range("I105").Value = 1818.9
$ws.Range("K105").Value = 1818.9
$ws.Range("M105").Value = -71.90000000000009
$ws.Range("H126").Value = 7333.8
$ws.Range("I126").Value = 7955.4546
$ws.Range("K126").Value = 23866.3638
$ws.Range("M126").Value = -21396.3638

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2142.7715
$ws.Range("I102").Value = 1375.44
$ws.Range("K102").Value = 1375.44
$ws.Range("M102").Value = 246.5599999999999
$ws.Range("H122").Value = 2352.3809
$ws.Range("I122").Value = 2182.1428
$ws.Range("K122").Value = 6546.428400000001
$ws.Range("M122").Value = -4096.428400000001

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 2495.625
$ws.Range("J61").Value = 3678.75
$ws.Range("L61").Value = 3678.75
$ws.Range("N61").Value = -4082.75
$ws.Range("H111").Value = 103462.336
$ws.Range("J111").Value = 103462.336
$ws.Range("L111").Value = 103462.336
$ws.Range("N111").Value = -111642.336
$ws.Range("H113").Value = 2495.625
$ws.Range("J113").Value = 3678.75
$ws.Range("L113").Value = 3678.75
$ws.Range("N113").Value = -8018.75
$ws.Range("H122").Value = 4086.0908
$ws.Range("I122").Value = 3686.7097
$ws.Range("J122").Value = 5038.4614
$ws.Range("K122").Value = 11060.1291
$ws.Range("L122").Value = 15115.3842
$ws.Range("M122").Value = -8610.1291
$ws.Range("N122").Value = -20015.3842
$ws.Range("H136").Value = 4552.067
$ws.Range("I136").Value = 2557.7188
$ws.Range("J136").Value = 9461.23
$ws.Range("K136").Value = 7673.1564
$ws.Range("L136").Value = 28383.69
$ws.Range("M136").Value = -5123.1564
$ws.Range("N136").Value = -33483.69
$ws.Range("H141").Value = 0
$ws.Range("J141").Value = 0
$ws.Range("L141").Value = 0
$ws.Range("N141").Value = $null

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H45").Value = 30000
$ws.Range("I45").Value = 0
$ws.Range("J45").Value = 30000
$ws.Range("K45").Value = 0
$ws.Range("L45").Value = 30000
$ws.Range("M45").Value = $null
$ws.Range("N45").Value = -30982
$ws.Range("H81").Value = 1687.375
$ws.Range("I81").Value = 1642.7142
$ws.Range("K81").Value = 3285.4284
$ws.Range("M81").Value = -2224.4284
$ws.Range("H84").Value = 1687.375
$ws.Range("I84").Value = 1642.7142
$ws.Range("K84").Value = 16427.142
$ws.Range("M84").Value = -11123.142
$ws.Range("H108").Value = 40000
$ws.Range("J108").Value = 40000
$ws.Range("L108").Value = 40000
$ws.Range("N108").Value = -47680
